$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A482").Value = 4
$ws.Range("B482").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C482").Value = 'Los Lagos'
$ws.Range("D482").Value = 44656
$ws.Range("E482").Value = 10
$ws.Range("F482").Value = 'Fruta'
$ws.Range("G482").Value = 100103
$ws.Range("H482").Value = 'Frutos de hueso (carozo)'
$ws.Range("I482").Value = 100103006
$ws.Range("J482").Value = 'Nectarín'
$ws.Range("K482").Value = 'Artic Snow'
$ws.Range("L482").Value = 'Especial'
$ws.Range("M482").Value = 200
$ws.Range("N482").Value = 24000
$ws.Range("O482").Value = 24000
$ws.Range("P482").Value = 24000
$ws.Range("Q482").Value = '$/caja 15 kilos empedrada'
$ws.Range("R482").Value = 'Región de O''Higgins'
$ws.Range("S482").Value = 1600
$ws.Range("T482").Value = 15

$ws.Range("A483").Value = 4
$ws.Range("B483").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C483").Value = 'Los Lagos'
$ws.Range("D483").Value = 44656
$ws.Range("E483").Value = 10
$ws.Range("F483").Value = 'Fruta'
$ws.Range("G483").Value = 100103
$ws.Range("H483").Value = 'Frutos de hueso (carozo)'
$ws.Range("I483").Value = 100103006
$ws.Range("J483").Value = 'Nectarín'
$ws.Range("K483").Value = 'Artic Snow'
$ws.Range("L483").Value = 'Primera'
$ws.Range("M483").Value = 200
$ws.Range("N483").Value = 20000
$ws.Range("O483").Value = 20000
$ws.Range("P483").Value = 20000
$ws.Range("Q483").Value = '$/caja 15 kilos empedrada'
$ws.Range("R483").Value = 'Región de O''Higgins'
$ws.Range("S483").Value = 1333
$ws.Range("T483").Value = 15

$ws.Range("A484").Value = 4
$ws.Range("B484").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C484").Value = 'Los Lagos'
$ws.Range("D484").Value = 44656
$ws.Range("E484").Value = 10
$ws.Range("F484").Value = 'Fruta'
$ws.Range("G484").Value = 100103
$ws.Range("H484").Value = 'Frutos de hueso (carozo)'
$ws.Range("I484").Value = 100103006
$ws.Range("J484").Value = 'Nectarín'
$ws.Range("K484").Value = 'Artic Snow'
$ws.Range("L484").Value = 'Segunda'
$ws.Range("M484").Value = 200
$ws.Range("N484").Value = 16000
$ws.Range("O484").Value = 16000
$ws.Range("P484").Value = 16000
$ws.Range("Q484").Value = '$/caja 15 kilos empedrada'
$ws.Range("R484").Value = 'Región de O''Higgins'
$ws.Range("S484").Value = 1067
$ws.Range("T484").Value = 15

$ws.Range("A485").Value = 4
$ws.Range("B485").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C485").Value = 'Los Lagos'
$ws.Range("D485").Value = 44656
$ws.Range("E485").Value = 10
$ws.Range("F485").Value = 'Fruta'
$ws.Range("G485").Value = 100103
$ws.Range("H485").Value = 'Frutos de hueso (carozo)'
$ws.Range("I485").Value = 100103006
$ws.Range("J485").Value = 'Nectarín'
$ws.Range("K485").Value = 'August Red'
$ws.Range("L485").Value = 'Especial'
$ws.Range("M485").Value = 200
$ws.Range("N485").Value = 24000
$ws.Range("O485").Value = 24000
$ws.Range("P485").Value = 24000
$ws.Range("Q485").Value = '$/caja 15 kilos empedrada'
$ws.Range("R485").Value = 'Región de O''Higgins'
$ws.Range("S485").Value = 1600
$ws.Range("T485").Value = 15

$ws.Range("A486").Value = 4
$ws.Range("B486").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C486").Value = 'Los Lagos'
$ws.Range("D486").Value = 44656
$ws.Range("E486").Value = 10
$ws.Range("F486").Value = 'Fruta'
$ws.Range("G486").Value = 100103
$ws.Range("H486").Value = 'Frutos de hueso (carozo)'
$ws.Range("I486").Value = 100103006
$ws.Range("J486").Value = 'Nectarín'
$ws.Range("K486").Value = 'August Red'
$ws.Range("L486").Value = 'Primera'
$ws.Range("M486").Value = 200
$ws.Range("N486").Value = 20000
$ws.Range("O486").Value = 20000
$ws.Range("P486").Value = 20000
$ws.Range("Q486").Value = '$/caja 15 kilos empedrada'
$ws.Range("R486").Value = 'Región de O''Higgins'
$ws.Range("S486").Value = 1333
$ws.Range("T486").Value = 15

$ws.Range("A487").Value = 4
$ws.Range("B487").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C487").Value = 'Los Lagos'
$ws.Range("D487").Value = 44656
$ws.Range("E487").Value = 10
$ws.Range("F487").Value = 'Fruta'
$ws.Range("G487").Value = 100103
$ws.Range("H487").Value = 'Frutos de hueso (carozo)'
$ws.Range("I487").Value = 100103006
$ws.Range("J487").Value = 'Nectarín'
$ws.Range("K487").Value = 'August Red'
$ws.Range("L487").Value = 'Segunda'
$ws.Range("M487").Value = 200
$ws.Range("N487").Value = 16000
$ws.Range("O487").Value = 16000
$ws.Range("P487").Value = 16000
$ws.Range("Q487").Value = '$/caja 15 kilos empedrada'
$ws.Range("R487").Value = 'Región de O''Higgins'
$ws.Range("S487").Value = 1067
$ws.Range("T487").Value = 15

$ws.Range("D482").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D483").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D484").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D485").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D486").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D487").NumberFormat = "YYYY-MM-DD HH:MM:SS"
